$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H28").Value = 571.875
$ws.Range("J28").Value = 577.6
$ws.Range("L28").Value = 577.6
$ws.Range("N28").Value = -1547.6

$ws.Range("H93").Value = 97437
$ws.Range("J93").Value = 97437
$ws.Range("L93").Value = 97437
$ws.Range("N93").Value = -102429

$ws.Range("H96").Value = 1025.3
$ws.Range("I96").Value = 769
$ws.Range("J96").Value = 3332
$ws.Range("K96").Value = 2307
$ws.Range("L96").Value = 9996
$ws.Range("M96").Value = -934
$ws.Range("N96").Value = -12742

$ws.Range("H98").Value = 250000720
$ws.Range("I98").Value = 500000450
$ws.Range("K98").Value = 500000450
$ws.Range("M98").Value = -499998952

$ws.Range("H100").Value = 4749.75
$ws.Range("I100").Value = 5000
$ws.Range("K100").Value = 5000
$ws.Range("M100").Value = -4459

$ws.Range("H106").Value = 24800.4
$ws.Range("I106").Value = 0
$ws.Range("J106").Value = 24800.4
$ws.Range("K106").Value = 0
$ws.Range("L106").Value = 24800.4
$ws.Range("N106").Value = -26062.4
$ws.Range("M106").ClearContents()

$ws.Range("H107").Value = 691.5
$ws.Range("I107").Value = 663.4545000000001
$ws.Range("K107").Value = 663.4545000000001
$ws.Range("M107").Value = 1256.5455

$ws.Range("H116").Value = 6064.5293
$ws.Range("I116").Value = 5495.5386
$ws.Range("K116").Value = 5495.5386
$ws.Range("M116").Value = -2053.5386

$ws.Range("H122").Value = 250000720
$ws.Range("I122").Value = 500000450
$ws.Range("K122").Value = 1500001350
$ws.Range("M122").Value = -1499998900

$ws.Range("H138").Value = 3424.5854
$ws.Range("I138").Value = 2636.6924
$ws.Range("J138").Value = 3790.3928
$ws.Range("K138").Value = 7910.0772
$ws.Range("L138").Value = 11371.1784
$ws.Range("M138").Value = -2770.0772
$ws.Range("N138").Value = -21651.1784

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H5").Value = 335
$ws.Range("J5").Value = 154
$ws.Range("L5").Value = 154
$ws.Range("N5").Value = -378

$ws.Range("H102").Value = 7399.625
$ws.Range("I102").Value = 7358.933
$ws.Range("K102").Value = 7358.933
$ws.Range("M102").Value = -5736.933

$ws.Range("H110").Value = 2118.5
$ws.Range("I110").Value = 2135.5715
$ws.Range("J110").Value = 1999
$ws.Range("K110").Value = 2135.5715
$ws.Range("L110").Value = 1999
$ws.Range("M110").Value = -90.57150000000001
$ws.Range("N110").Value = -6089

$ws.Range("H122").Value = 3974.9412
$ws.Range("J122").Value = 4536.615
$ws.Range("L122").Value = 13609.845
$ws.Range("N122").Value = -18509.845

$ws.Range("H141").Value = 119999
$ws.Range("J141").Value = 119999
$ws.Range("L141").Value = 119999
$ws.Range("N141").Value = -130359

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H4").Value = 335
$ws.Range("J4").Value = 154
$ws.Range("L4").Value = 154
$ws.Range("N4").Value = -384

$ws.Range("H20").Value = 2250
$ws.Range("I20").Value = 2000
$ws.Range("K20").Value = 2000
$ws.Range("M20").Value = -1753

$ws.Range("H105").Value = 0
$ws.Range("I105").Value = 0
$ws.Range("K105").Value = 0
$ws.Range("M105").ClearContents()

$ws.Range("H107").Value = 1910.5625
$ws.Range("I107").Value = 1762.0714
$ws.Range("K107").Value = 1762.0714
$ws.Range("M107").Value = 157.9286

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 19579.8
$ws.Range("I16").Value = 19579.8
$ws.Range("K16").Value = 19579.8
$ws.Range("M16").Value = -19292.8

$ws.Range("H31").Value = 872964.25
$ws.Range("I31").Value = 15060.546
$ws.Range("K31").Value = 15060.546
$ws.Range("M31").Value = -14765.546

$ws.Range("H34").Value = 872964.25
$ws.Range("I34").Value = 15060.546
$ws.Range("K34").Value = 15060.546
$ws.Range("M34").Value = -14858.546

$ws.Range("H39").Value = 12728.714
$ws.Range("I39").Value = 3033.6667
$ws.Range("K39").Value = 3033.6667
$ws.Range("M39").Value = -2642.6667

$ws.Range("H49").Value = 12728.714
$ws.Range("I49").Value = 3033.6667
$ws.Range("K49").Value = 3033.6667
$ws.Range("M49").Value = -2851.6667

$ws.Range("H113").Value = 19579.8
$ws.Range("I113").Value = 19579.8
$ws.Range("K113").Value = 19579.8
$ws.Range("M113").Value = -17409.8

$ws.Range("H122").Value = 3114.182
$ws.Range("I122").Value = 1933.3334
$ws.Range("K122").Value = 5800.0002
$ws.Range("M122").Value = -3350.0002

$ws.Range("H132").Value = 4153.5186
$ws.Range("J132").Value = 3900
$ws.Range("L132").Value = 11700
$ws.Range("N132").Value = -16760

$ws.Range("H135").Value = 0
$ws.Range("J135").Value = 0
$ws.Range("L135").Value = 0
$ws.Range("N135").ClearContents()

$ws.Range("H141").Value = 94979.27
$ws.Range("J141").Value = 106086
$ws.Range("L141").Value = 106086
$ws.Range("N141").Value = -116446

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H34").Value = 1599.6923
$ws.Range("I34").Value = 311.5
$ws.Range("J34").Value = 2172.2222
$ws.Range("K34").Value = 934.5
$ws.Range("L34").Value = 6516.6666
$ws.Range("M34").Value = -850.5
$ws.Range("N34").Value = -6684.6666

$ws.Range("H39").Value = 110894.9
$ws.Range("J39").Value = 213799.8
$ws.Range("L39").Value = 641399.3999999999
$ws.Range("N39").Value = -641987.3999999999

$ws.Range("H56").Value = 7661.9473
$ws.Range("I56").Value = 7661.9473
$ws.Range("K56").Value = 7661.9473
$ws.Range("M56").Value = -7131.9473

$ws.Range("H96").Value = 30000
$ws.Range("J96").Value = 30000
$ws.Range("L96").Value = 90000
$ws.Range("N96").Value = -94118

$ws.Range("H113").Value = 1308.28
$ws.Range("I113").Value = 724.625
$ws.Range("J113").Value = 1582.9412
$ws.Range("K113").Value = 2173.875
$ws.Range("L113").Value = 4748.8236
$ws.Range("M113").Value = -3.875
$ws.Range("N113").Value = -9088.8236

$ws.Range("H122").Value = 631.5
$ws.Range("J122").Value = 656.0769
$ws.Range("L122").Value = 5904.6921
$ws.Range("N122").Value = -10804.6921

$ws.Range("H125").Value = 15505.5
$ws.Range("J125").Value = 15505.5
$ws.Range("L125").Value = 46516.5
$ws.Range("N125").Value = -56356.5

$ws.Range("H131").Value = 1726.6923
$ws.Range("I131").Value = 1794
$ws.Range("J131").Value = 1619
$ws.Range("K131").Value = 5382
$ws.Range("L131").Value = 4857
$ws.Range("M131").Value = -342
$ws.Range("N131").Value = -14937

$ws.Range("H132").Value = 2346.88
$ws.Range("I132").Value = 2417.5386
$ws.Range("J132").Value = 2270.3333
$ws.Range("K132").Value = 21757.8474
$ws.Range("L132").Value = 20432.9997
$ws.Range("M132").Value = -19227.8474
$ws.Range("N132").Value = -25492.9997

$ws.Range("H141").Value = 6926.727
$ws.Range("I141").Value = 4068.913
$ws.Range("K141").Value = 12206.739
$ws.Range("M141").Value = -7026.739

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H43").Value = 17254.25
$ws.Range("I43").Value = 9672.333000000001
$ws.Range("K43").Value = 9672.333000000001
$ws.Range("M43").Value = -9521.333000000001

$ws.Range("H44").Value = 39999.75
$ws.Range("I44").Value = 39999
$ws.Range("K44").Value = 39999
$ws.Range("M44").Value = -39403

$ws.Range("H52").Value = 0
$ws.Range("J52").Value = 0
$ws.Range("L52").Value = 0
$ws.Range("N52").ClearContents()

$ws.Range("H55").Value = 10000
$ws.Range("I55").Value = 10000
$ws.Range("K55").Value = 10000
$ws.Range("M55").Value = -9673

$ws.Range("H132").Value = 35716860
$ws.Range("I132").Value = 35716860
$ws.Range("K132").Value = 107150580
$ws.Range("M132").Value = -107148050

$ws.Range("H135").Value = 194663
$ws.Range("J135").Value = 194663
$ws.Range("L135").Value = 194663
$ws.Range("N135").Value = -204803

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4236841.5
$ws.Range("I132").Value = 505903.56
$ws.Range("K132").Value = 1517710.68
$ws.Range("M132").Value = -1515180.68

$ws.Range("H136").Value = 110418.38
$ws.Range("I136").Value = 81368.16
$ws.Range("K136").Value = 244104.48
$ws.Range("M136").Value = -241554.48

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 37782.5
$ws.Range("I54").Value = 37070
$ws.Range("K54").Value = 37070
$ws.Range("M54").Value = -36550

$ws.Range("H113").Value = 700.27026
$ws.Range("I113").Value = 627.88464
$ws.Range("J113").Value = 871.36365
$ws.Range("K113").Value = 1883.65392
$ws.Range("L113").Value = 2614.09095
$ws.Range("M113").Value = 286.34608
$ws.Range("N113").Value = -6954.09095

$ws.Range("H117").Value = 111000
$ws.Range("J117").Value = 111000
$ws.Range("L117").Value = 111000
$ws.Range("N117").Value = -120178

$ws.Range("H126").Value = 4217.457
$ws.Range("I126").Value = 3262.762
$ws.Range("K126").Value = 9788.286
$ws.Range("M126").Value = -7318.286

$ws.Range("H132").Value = 306887.94
$ws.Range("I132").Value = 2604.1738
$ws.Range("K132").Value = 7812.5214
$ws.Range("M132").Value = -5282.5214

$ws.Range("H136").Value = 4050.8462
$ws.Range("I136").Value = 1287.2858
$ws.Range("K136").Value = 3861.8574
$ws.Range("M136").Value = -1311.8574
